$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Column layout: before had cols 1-2 sharing one <col> element
#    (min="1" max="2"), after the edit column A gets its own <col> entry.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 30.7109375

# ---------------------------------------------------------------------------
# 2) Insert a new row above row 13. This shifts the old rows 13-21 down to
#    14-22 (carrying their row heights and cell styles with them), which
#    already lines every column-A label up with its correct row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# ---------------------------------------------------------------------------
# 3) Fill in the new / corrected long-form text for columns B and C.
# ---------------------------------------------------------------------------
$objetivos = @'
Apresentar ao aluno as aplicações dos geossintéticos em obras de proteção e recuperação ambiental dando ênfase às múltiplas funções dos geossintéticos. Aspectos como a drenagem e filtração, a separação, barreiras de silte, reforço de solos e os sistemas de contenção de resíduos e de efluentes que podem ser utilizados com vantagens técnicas e econômicas tanto em obras de proteção ambiental como na recuperação de áreas degradadas e/ou contaminadas.
'@

$resumido = @'
Obras de proteção ambiental: sistemas de contenção de resíduos; lagoas de efluentes, remediação e mitigação de áreas degradadas; Geossintéticos: tipos e aplicações; Geotêxteis, geomantas e geocompostos para a drenagem: tipos, propriedades, ensaios caracterização e desempenho; Drenagem e filtração com geossintéticos. Critérios de filtração com geossintéticos. Detalhes construtivos de obras de drenagem com geossintéticos; Sistemas de drenagem em obras de proteção e de recuperação ambiental; Geogrelhas e geotêxteis para reforço: tipos, propriedades, ensaios de caracterização e de desempenho; Princípios gerais de reforço de solos. Reforço de solos com geossintéticos. Estruturas em solos reforçados. Detalhamento da construção de obras em solo reforçado; Geomembranas: tipos, propriedades, ensaios de caracterização e de desempenho; Geocomposto bentonítico: tipos, propriedades, ensaios de caracterização e de desempenho; Barreiras impermeabilizantes de fundação e de cobertura.
'@

$completo = @'
Obras de proteção ambiental: sistemas de contenção de resíduos; lagoas de efluentes, remediação e mitigação de áreas degradadas; Geossintéticos: tipos e aplicações; Geotéxteis, geomantas e geocompostos para a drenagem: tipos, propriedades; ensaios de caracterização e de desempenho; Drenagem e filtração com geossintéticos. Critérios de filtração com geossintéticos. Detalhes construtivos de obras de drenagem com geossintéticos; Sistemas de drenagem em obras de proteção e de recuperação ambiental; Geogrelhas e geotêxteis para reforço; tipos, propriedades, ensaios de caracterização e de desempenho; Princípios gerais de reforço de solos. Reforço de solos com geossintéticos. Estruturas em solos reforçados. Detalhamento da construção de obras em solo reforçado; Geomembranas: tipos, propriedades, ensaios de caracterização e de desempenho; Geocomposto bentonítico: tipos, propriedades, ensaios de caracterização e de desempenho; Barreiras impermeabilizantes de fundação e de cobertura: tipos, especificações, características; Solicitações físicas mecânicas e químicas de geomembranas em sistemas de contenção de resíduos e de efluentes; Aspectos construtivos de barreiras impermeabilizantes e de cobertura; especificações de geossintéticos.
'@

$biblio = @'
SHARMA, H. D. & LEWIS, S.P. (1994) Waste Containement Systems, Waste Stabilization and Landfils. Design and Evaluation - Joh Willy & Sons, Inc. New York;
KOERNER. R.M. (1997) Designing with Geosynthetics, Prentice Hall Inc.;
QIAN, X/; KOERNER, R.M. & GRAY, D.H. (2002) Geotechnical Aspects of Landfill Design and Construction - Prentice Hall, Upper Saddle Riner, New Jersey 07458;
Manual Brasileiro de Geossintéticos, ABINT, 2004;
INGOLD, T.S. (1994) The Geotêxtiles and Geomembranes. Manual, Elsevier, London, 610p.;
VAN ZATEN, R.V. (1986) Geotextiles and Geomembranes in Civil Engineering. A.A. Balkema, Rotherdam, netherlands, 654p.
'@

$aulas = "Aulas expositivas; microcomputadores; seminários; visitas técnicas."
$media = "Média ponderada de notas de provas e seminários."
$prova = "Prova única com nota igual ou superior a 5,0 (cinco)."
$docente = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Row 10 - Objetivos: long description replaces the old (misplaced) docente text
$ws.Range("B10:C10").Value = $objetivos

# Row 13 (new, blank row) - only B/C hold the docente name, no label in A
$ws.Range("B13:C13").Value = $docente

# Row 14 - Programa resumido: replaces "Semestral" placeholder
$ws.Range("B14:C14").Value = $resumido

# Row 16 - Programa: replaces misplaced date text
$ws.Range("B16:C16").Value = $completo

# Row 19 - Método: now holds the "Aulas expositivas..." text
$ws.Range("B19:C19").Value = $aulas

# Row 20 - Critério: now holds the "Média ponderada..." text
$ws.Range("B20:C20").Value = $media

# Row 21 - Norma de recuperação: now holds the "Prova única..." text
$ws.Range("B21:C21").Value = $prova

# Row 22 - Bibliografia: new bibliography text (label already shifted into place)
$ws.Range("B22:C22").Value = $biblio
